$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2019704433497537
$ws.Range("C2").Value = 0.5566502463054187
$ws.Range("J2").Value = 0.01477832512315271
$ws.Range("P2").Value = 0.1527093596059113
$ws.Range("S2").Value = 0.07389162561576355
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1785714285714286
$ws.Range("B6").Value = 0.09268292682926829
$ws.Range("D6").Value = 0.01463414634146342
$ws.Range("F6").Value = 0.03902439024390244
$ws.Range("J6").Value = 0.2292682926829268
$ws.Range("O6").Value = 0.03414634146341464
$ws.Range("Q6").Value = 0.1560975609756098
$ws.Range("R6").Value = 0.07317073170731707
$ws.Range("S6").Value = 0.3609756097560975
$ws.Range("B7").Value = 0.1029411764705882
$ws.Range("D7").Value = 0.01470588235294118
$ws.Range("F7").Value = 0.07352941176470588
$ws.Range("J7").Value = 0.1102941176470588
$ws.Range("O7").Value = 0.01470588235294118
$ws.Range("Q7").Value = 0.1617647058823529
$ws.Range("R7").Value = 0.1323529411764706
$ws.Range("S7").Value = 0.3897058823529412
$ws.Range("B8").Value = 0.05847953216374269
$ws.Range("D8").Value = 0.01754385964912281
$ws.Range("F8").Value = 0.04678362573099415
$ws.Range("J8").Value = 0.108187134502924
$ws.Range("O8").Value = 0.01461988304093567
$ws.Range("Q8").Value = 0.2192982456140351
$ws.Range("R8").Value = 0.1169590643274854
$ws.Range("S8").Value = 0.4181286549707602
$ws.Range("B9").Value = 0.06779661016949153
$ws.Range("D9").Value = 0.005649717514124294
$ws.Range("F9").Value = 0.06214689265536723
$ws.Range("J9").Value = 0.1016949152542373
$ws.Range("O9").Value = 0.03389830508474576
$ws.Range("Q9").Value = 0.1468926553672316
$ws.Range("R9").Value = 0.1355932203389831
$ws.Range("S9").Value = 0.4463276836158192
$ws.Range("B10").Value = 0.07933884297520662
$ws.Range("D10").Value = 0.01570247933884298
$ws.Range("E10").Value = 0.0008264462809917355
$ws.Range("F10").Value = 0.07768595041322314
$ws.Range("J10").Value = 0.09669421487603305
$ws.Range("O10").Value = 0.01983471074380165
$ws.Range("Q10").Value = 0.2247933884297521
$ws.Range("R10").Value = 0.1148760330578512
$ws.Range("S10").Value = 0.3702479338842975
$ws.Range("G11").Value = 0.1205357142857143
$ws.Range("J11").Value = 0.1071428571428571
$ws.Range("K11").Value = 0.1919642857142857
$ws.Range("L11").Value = 0.5401785714285714
$ws.Range("S11").Value = 0.04017857142857143
$ws.Range("G12").Value = 0.734375
$ws.Range("J12").Value = 0.171875
$ws.Range("K12").Value = 0.0078125
$ws.Range("L12").Value = 0.0390625
$ws.Range("S12").Value = 0.046875
$ws.Range("F15").Value = 0.04054054054054054
$ws.Range("H15").Value = 0.1756756756756757
$ws.Range("I15").Value = 0.08108108108108109
$ws.Range("J15").Value = 0.3783783783783784
$ws.Range("K15").Value = 0.04504504504504504
$ws.Range("M15").Value = 0.004504504504504504
$ws.Range("O15").Value = 0.04504504504504504
$ws.Range("S15").Value = 0.2297297297297297
$ws.Range("F16").Value = 0.02419354838709677
$ws.Range("H16").Value = 0.06451612903225806
$ws.Range("I16").Value = 0.04032258064516129
$ws.Range("J16").Value = 0.5887096774193549
$ws.Range("K16").Value = 0.08870967741935484
$ws.Range("M16").Value = 0.03225806451612903
$ws.Range("N16").Value = 0.008064516129032258
$ws.Range("O16").Value = 0.08870967741935484
$ws.Range("S16").Value = 0.06451612903225806
$ws.Range("F17").Value = 0.0188235294117647
$ws.Range("H17").Value = 0.1388235294117647
$ws.Range("I17").Value = 0.09176470588235294
$ws.Range("J17").Value = 0.48
$ws.Range("K17").Value = 0.09647058823529411
$ws.Range("M17").Value = 0.02352941176470588
$ws.Range("N17").Value = 0.002352941176470588
$ws.Range("O17").Value = 0.06352941176470588
$ws.Range("S17").Value = 0.08470588235294117
$ws.Range("F18").Value = 0.008547008547008548
$ws.Range("H18").Value = 0.1239316239316239
$ws.Range("I18").Value = 0.09829059829059829
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.0811965811965812
$ws.Range("M18").Value = 0.02136752136752137
$ws.Range("O18").Value = 0.0811965811965812
$ws.Range("S18").Value = 0.08547008547008547
$ws.Range("F19").Value = 0.01424501424501425
$ws.Range("H19").Value = 0.2003798670465337
$ws.Range("I19").Value = 0.08736942070275404
$ws.Range("J19").Value = 0.4235517568850902
$ws.Range("K19").Value = 0.09211775878442545
$ws.Range("M19").Value = 0.0113960113960114
$ws.Range("N19").Value = 0.001899335232668566
$ws.Range("O19").Value = 0.08452041785375118
$ws.Range("S19").Value = 0.08452041785375118
